# Rename the "_old" / "_new" header-column suffixes to the format-version
# specific suffixes ("_FV2410" / "_FV2504") used by this export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace "_old$", "_FV2410")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace "_new$", "_FV2504")
}

# Turn the used range into an actual Excel Table ("Table1") so the header
# row gets AutoFilter dropdowns and the column names become structured
# table-column names.
$usedRange = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $usedRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit complete"
